# Update column G ("K") values on Sheet1 per the regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gValues = @{
    "2" = 1
    "3" = 1
    "4" = 1
    "5" = 2
    "6" = 1
    "7" = 2
    "8" = 0
    "9" = 0
    "10" = 0
    "11" = 1
    "12" = 2
    "13" = 1
    "14" = 0
    "15" = 3
    "16" = 0
    "17" = 0
    "18" = 0
    "19" = 0
    "20" = 0
    "21" = 0
    "22" = 2
    "24" = 0
    "25" = 0
    "26" = 2
    "27" = 1
    "28" = 0
    "29" = 2
    "30" = 1
    "31" = 0
    "32" = 0
    "33" = 1
    "34" = 3
    "35" = 1
    "36" = 2
    "37" = 0
    "38" = 3
    "39" = 1
    "40" = 1
    "41" = 1
    "43" = 1
    "44" = 1
    "45" = 1
    "46" = 2
    "47" = 3
    "48" = 0
    "49" = 1
    "50" = 0
    "51" = 0
    "52" = 1
    "53" = 1
    "54" = 0
    "55" = 3
    "56" = 2
    "57" = 3
    "58" = 0
    "59" = 1
    "60" = 1
    "61" = 1
    "62" = 1
    "63" = 0
    "64" = 0
    "65" = 1
    "66" = 1
    "67" = 3
    "68" = 1
    "70" = 1
    "72" = 2
    "73" = 1
    "74" = 2
    "75" = 1
    "76" = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
